$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dscanbo")

$data = @(
    @("00123456", "Nguyễn Ân Hiên",   "Hệ thống thông tin",                "Công nghệ thông tin và truyền thông", "vana@gmail.com"),
    @("00123457", "Trần Văn Bình",    "Công nghệ thông tin",               "Công nghệ thông tin và truyền thông", "vanb@gmail.com"),
    @("00123458", "Phan Công Huy",    "Mạng máy tính và TT",               "Công nghệ thông tin và truyền thông", "pchuy@gmail.com"),
    @("00123459", "Lê Thị Hậu",       "Khoa học máy tính",                 "Công nghệ thông tin và truyền thông", "hau@hotmail.com"),
    @("00123460", "Nguyễn Chín",      "Mạng máy tính và TT",               "Công nghệ thông tin và truyền thông", "chin@email.com"),
    @("00123461", "Lê Thành Công",    "Công nghệ thông tin",               "Công nghệ thông tin và truyền thông", "cong@outlook.com"),
    @("00123462", "Phương Tín",       "Tin học ứng dụng",                  "Công nghệ thông tin và truyền thông", "Ptin@gmail.com.vn"),
    @("00123463", "Lâm Định Cương",   "Tin học ứng dụng",                  "Công nghệ thông tin và truyền thông", "ldcuong@gmail.com"),
    @("00123464", "Mai Ất Nguyên",    "Công nghệ phần mềm",                "Công nghệ thông tin và truyền thông", "nguyenmaiat@gmail.com"),
    @("00123465", "Phạm Thắng",       "Mạng máy tính và TT",               "Công nghệ thông tin và truyền thông", "thangthang@gmail.com.vn"),
    @("00123466", "Trương Tuấn Kiên", "Hệ thống thông tin",                "Công nghệ thông tin và truyền thông", "tkien@yahoo.com.vn"),
    @("00123467", "Đinh Phùng Quân",  "Hệ thống thông tin",                "Công nghệ thông tin và truyền thông", "phungquan@gmail.com.vn")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $row++
}

$ws.Range("I20").Select()
